# SwaadSutra_Consolidated_2026-01-20.xlsx update
# New order #22 (Pooja, Til Poli x1) came in at 2026-01-20 11:13 and is
# inserted as the newest row at the top of the "All Orders" log (row 2),
# pushing all existing order rows down by one. The "Daily Summary" sheet's
# 2026-01-20 totals are updated to reflect the new order.

$wb = $excel.ActiveWorkbook

# ---- Sheet: All Orders ----
$ws = $wb.Worksheets.Item("All Orders")

# Insert a new blank row at row 2, shifting existing order rows (and their
# formatting) down by one - row 2 becomes row 3, row 3 becomes row 4, etc.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the incoming order's data.
$ws.Cells.Item(2, 1).Value = 22                    # Order ID (numeric)
$ws.Cells.Item(2, 2).Value = "2026-01-20 11:13"    # Date
$ws.Cells.Item(2, 3).Value = "Pooja"               # Customer

# Flat No / Phone / Collection Date look like numbers or dates, so force
# text formatting first to keep them stored as text (matching the rest of
# the sheet, e.g. phone numbers elsewhere in the log).
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "12"                  # Flat No

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "9096648553"          # Phone

$ws.Cells.Item(2, 6).Value = "Til Poli x1"         # Items
$ws.Cells.Item(2, 7).Value = 30                    # Total (numeric)
$ws.Cells.Item(2, 8).Value = "NEW"                 # Status
$ws.Cells.Item(2, 9).Value = "PENDING"             # Payment

$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2026-01-20"         # Collection Date

$ws.Cells.Item(2, 11).Value = "16:43"              # Collection Time
$ws.Cells.Item(2, 12).Value = ""                   # Notes
$ws.Cells.Item(2, 13).Value = ""                   # Cancel Reason
$ws.Cells.Item(2, 14).Value = ""                   # Feedback

# ---- Sheet: Daily Summary ----
$ws2 = $wb.Worksheets.Item("Daily Summary")

# Row 2 corresponds to 2026-01-20; add the new NEW/PENDING order of 30.
$ws2.Cells.Item(2, 2).Value = 3                    # Total Orders: 2 -> 3
$ws2.Cells.Item(2, 5).Value = 80                   # Revenue: 50 -> 80
$ws2.Cells.Item(2, 7).Value = 30                   # Pending: 0 -> 30

Write-Host "Applied SwaadSutra 2026-01-20 11:13 update"
